# Added filtering options for the Component Analysis
# Clear specific cells that are no longer part of the filtered
# Component Analysis output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @("J2", "K2", "I3", "J3", "K3", "K5", "J6", "K6", "I7", "J7", "K7")

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
